# Populate a "Date" key/value pair on "Test Sheet 2" (rows 10-11) and make
# that sheet the active one (mirrors the author's manual edit + tab switch).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Sheet 2")

# New key label in C10.
$ws.Range("C10").Value = "Date"

# D10: a date formatted as DD/MM/YY (serial 51920 -> 23/02/42).
$ws.Range("D10").Value = 51920
$ws.Range("D10").NumberFormat = "DD/MM/YY"

# D11: a date formatted as "D MMM YYYY" (serial 43466 -> 1 Jan 2019).
$ws.Range("D11").Value = 43466
$ws.Range("D11").NumberFormat = "D\ MMM\ YYYY"

# Switch the active sheet/selection to Test Sheet 2, cell D11 (last edited).
$ws.Activate() | Out-Null
$ws.Range("D11").Select() | Out-Null
